$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after Sheet1 and name it Sheet2
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 with the login data (fill column A first, then column B, so the
# shared-strings table is built in the same order as the target workbook)
$ws2.Range("A1").Value = "用户名"
$ws2.Range("A2").Value = "密码"
$ws2.Range("B1").Value = "admin"
$ws2.Range("B2").Value = "admin123"

# Update selections: Sheet1 now has A1:A2 selected, Sheet2 has B3 selected (and is the active tab)
$ws1.Range("A1:A2").Select() | Out-Null
$ws2.Range("B3").Select() | Out-Null
$ws2.Activate() | Out-Null
